$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the 2 Ohm resistor value part number (row 18, column B)
$ws.Range("B18").Value = "TBH25P2R00JE"

# Add new rows for heatsink and fan mounting parts (value-only, no reference)
$ws.Range("B26").Value = "Heatsink: 542502B00000G x16"
$ws.Range("B27").Value = "Fanx: 2x 0D7025-24MB x2"

# Update the active selection to match the final state
$ws.Range("E18").Select()
